$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45932, 4000.97075462538, 5595.47713359752, 6020, 5016.860973, 24.6403063321726),
    @(45933, 4040.15146752679, 5189.39918975202, 3620, 5067.893254, 108.214207342718),
    @(45934, 860.401107928437, 3160.29772488092, 3620, 1664.9207, 14.3673882063533),
    @(45935, 747.320273491827, 2804.99207822295, 3620, 1528.13076, -1.42489313620304),
    @(45936, 4031.01198509058, 5241.15192381869, 3620, 5046.137519, 109.844894072005),
    @(45937, 4367.23135895568, 5185.48425572587, 3620, 6288.409189, 145.277586907091),
    @(45938, 4367.23135895568, 5094.35958788455, 3620, 6288.409189, 141.480725747036),
    @(45939, 4367.23135895568, 4905.23510647766, 3620, 6288.409189, 133.600539021749),
    @(45940, 5468.88248205865, 4241.83480403972, 3620, 8246.898597, 141.660454957545),
    @(45941, 1519.17029469927, 2172.63364585955, 3620, 3912.04371, 39.396127548345),
    @(45942, 1415.02104924416, 2214.47941371422, 3620, 3798.778212, 40.7598573529191),
    @(45943, 5868.78419466605, 5096.4664700444, 3620, 9003.232851, 192.121463599098),
    @(45944, 5868.78419466605, 5513.0953049223, 3620, 9003.232851, 209.480998385677),
    @(45945, 5868.78419466605, 5691.10258433736, 3620, 9024.242837, 217.773384444638)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
